$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Filtros Descontinuados"
$ws.Range("B3").Value = "Discontinued Filters"
$ws.Range("C3").Value = "Categoría de filtros que ya no se fabrican"
$ws.Range("D3").Value = "Category of filters that are no longer manufactured"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "inactive"
